# Update the "Raw" worksheet per the 2022-end-of-year status refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw")

# ccao-condominium.R was renamed to ccao-condominium_parking.R
$ws.Range("A3").Value = "ccao-condominium_parking.R"

# New "Issue Preventing Successful Run" notes
$ws.Range("D2").Value = "waitin on valuations"
$ws.Range("D3").Value = "no need to update, 12/8/22"
$ws.Range("D10").Value = "waiting on great schools api"
$ws.Range("D20").Value = "waiting to hear from Ray on subdivisions"
$ws.Range("D21").Value = "waiting to hear from Josh on 2022 parcels"
$ws.Range("D22").Value = "nothing to update, 12/14/22"
$ws.Range("D24").Value = "waiting to hear from Josh on 21/22 school tax districts"
$ws.Range("D25").Value = "waiting to hear from Josh on 21/22 tax districts"

# spatial-transit.R now has a successful run logged, with a new note
# (copy the date formatting already used on other "Last Date Successfully Run" cells)
$ws.Range("B4").Copy($ws.Range("B26"))
$ws.Range("B26").Value = 44909
$ws.Range("C26").Value = "WRR"
$ws.Range("D26").Value = "PACE feed no longer maintained"

# Move the active selection to D10, matching where the author was last editing
$ws.Range("D10").Select()
